$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "Agile Co-Development Services"
$ws.Range("B2").Value = "GVT000ABC1234"
$ws.Range("B3").Value = "1 May 24 - 30"
$ws.Range("B7").Value = "DevOps Engineer"
$ws.Range("B48").Value = "09 - February - 2025"
